$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates ---
# A8: "Volume 30   Number  32" -> "Volume 30   Number  33"
$ws.Range("A8").Characters(21,2).Text = "33"

# C9: "Report Covering the Week  8/7/2023  Through  8/13/2023"
#     -> "Report Covering the Week  8/14/2023  Through  8/20/2023"
$ws.Range("C9").Characters(27,8).Text = "8/14/2023"
$ws.Range("C9").Characters(47,9).Text = "8/20/2023"

# --- Weekly crime statistics table updates (rows 14-30) ---
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = -33.333333333333
$ws.Range("F14").Value = 21
$ws.Range("G14").Value = 26
$ws.Range("H14").Value = -19.230769230769
$ws.Range("I14").Value = 252
$ws.Range("J14").Value = 283
$ws.Range("K14").Value = -10.95406360424
$ws.Range("L14").Value = -17.105263157894
$ws.Range("M14").Value = -27.586206896551
$ws.Range("N14").Value = -79.595141700404

$ws.Range("C15").Value = 28
$ws.Range("D15").Value = 26
$ws.Range("E15").Value = 7.692307692307
$ws.Range("F15").Value = 104
$ws.Range("G15").Value = 125
$ws.Range("H15").Value = -16.8
$ws.Range("I15").Value = 926
$ws.Range("J15").Value = 1051
$ws.Range("K15").Value = -11.893434823977
$ws.Range("L15").Value = -3.941908713692
$ws.Range("M15").Value = 9.585798816568
$ws.Range("N15").Value = -55.714968914395

$ws.Range("C16").Value = 351
$ws.Range("D16").Value = 370
$ws.Range("E16").Value = -5.135135135135
$ws.Range("F16").Value = 1405
$ws.Range("G16").Value = 1531
$ws.Range("H16").Value = -8.229915088177
$ws.Range("I16").Value = 10333
$ws.Range("J16").Value = 11023
$ws.Range("K16").Value = -6.259638936768
$ws.Range("L16").Value = 31.246030737965
$ws.Range("M16").Value = -11.608212147134
$ws.Range("N16").Value = -80.591295854542

$ws.Range("C17").Value = 563
$ws.Range("D17").Value = 503
$ws.Range("E17").Value = 11.928429423459
$ws.Range("F17").Value = 2301
$ws.Range("G17").Value = 2230
$ws.Range("H17").Value = 3.183856502242
$ws.Range("I17").Value = 17741
$ws.Range("J17").Value = 16793
$ws.Range("K17").Value = 5.645209313404
$ws.Range("L17").Value = 26.748588983353
$ws.Range("M17").Value = 61.41388408698
$ws.Range("N17").Value = -33.733004631704

$ws.Range("C18").Value = 265
$ws.Range("D18").Value = 301
$ws.Range("E18").Value = -11.960132890365
$ws.Range("F18").Value = 1098
$ws.Range("G18").Value = 1187
$ws.Range("H18").Value = -7.497893850042
$ws.Range("I18").Value = 8837
$ws.Range("J18").Value = 9853
$ws.Range("K18").Value = -10.311580229371
$ws.Range("L18").Value = 19.483504597079
$ws.Range("M18").Value = -22.814219582496
$ws.Range("N18").Value = -85.999904944471

$ws.Range("C19").Value = 1047
$ws.Range("D19").Value = 1098
$ws.Range("E19").Value = -4.644808743169
$ws.Range("F19").Value = 4186
$ws.Range("G19").Value = 4364
$ws.Range("H19").Value = -4.078826764436
$ws.Range("I19").Value = 31650
$ws.Range("J19").Value = 32410
$ws.Range("K19").Value = -2.344955260722
$ws.Range("L19").Value = 44.349174496032
$ws.Range("M19").Value = 36.422413793103
$ws.Range("N19").Value = -40.984523587544

$ws.Range("C20").Value = 338
$ws.Range("D20").Value = 228
$ws.Range("E20").Value = 48.245614035087
$ws.Range("F20").Value = 1369
$ws.Range("G20").Value = 1052
$ws.Range("H20").Value = 30.133079847908
$ws.Range("I20").Value = 9986
$ws.Range("J20").Value = 8372
$ws.Range("K20").Value = 19.278547539417
$ws.Range("L20").Value = 66.294754371357
$ws.Range("M20").Value = 52.995250497931
$ws.Range("N20").Value = -85.800213295414

$ws.Range("C21").Value = 2594
$ws.Range("D21").Value = 2529
$ws.Range("E21").Value = 2.570185844207
$ws.Range("F21").Value = 10484
$ws.Range("G21").Value = 10515
$ws.Range("H21").Value = -0.294816928197
$ws.Range("I21").Value = 79725
$ws.Range("J21").Value = 79785
$ws.Range("K21").Value = -0.075202105658
$ws.Range("L21").Value = 36.363636363636
$ws.Range("M21").Value = 22.559569561875
$ws.Range("N21").Value = -70.517319803411

$ws.Range("C22").Value = 42
$ws.Range("D22").Value = 42
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 163
$ws.Range("G22").Value = 150
$ws.Range("H22").Value = 8.666666666666
$ws.Range("I22").Value = 1385
$ws.Range("J22").Value = 1442
$ws.Range("K22").Value = -3.952843273231
$ws.Range("L22").Value = 43.078512396694
$ws.Range("M22").Value = 4.135338345864

$ws.Range("C23").Value = 109
$ws.Range("D23").Value = 109
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 495
$ws.Range("G23").Value = 512
$ws.Range("H23").Value = -3.3203125
$ws.Range("I23").Value = 3971
$ws.Range("J23").Value = 3846
$ws.Range("K23").Value = 3.2501300052
$ws.Range("L23").Value = 17.035072207486
$ws.Range("M23").Value = 51.564885496183

$ws.Range("C24").Value = 2211
$ws.Range("D24").Value = 2446
$ws.Range("E24").Value = -9.60752248569
$ws.Range("F24").Value = 8941
$ws.Range("G24").Value = 9657
$ws.Range("H24").Value = -7.414310862586
$ws.Range("I24").Value = 70017
$ws.Range("J24").Value = 72229
$ws.Range("K24").Value = -3.062481828628
$ws.Range("L24").Value = 38.160543036426
$ws.Range("M24").Value = 37.239797718452

$ws.Range("C25").Value = 857
$ws.Range("D25").Value = 726
$ws.Range("E25").Value = 18.044077134986
$ws.Range("F25").Value = 3558
$ws.Range("G25").Value = 3189
$ws.Range("H25").Value = 11.571025399811
$ws.Range("I25").Value = 27895
$ws.Range("J25").Value = 26543
$ws.Range("K25").Value = 5.093621670496
$ws.Range("L25").Value = 28.040943725328
$ws.Range("M25").Value = -6.029981472123

$ws.Range("C26").Value = 54
$ws.Range("D26").Value = 45
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 181
$ws.Range("G26").Value = 187
$ws.Range("H26").Value = -3.208556149732
$ws.Range("I26").Value = 1552
$ws.Range("J26").Value = 1692
$ws.Range("K26").Value = -8.274231678487
$ws.Range("L26").Value = -0.51282051282

$ws.Range("C27").Value = 105
$ws.Range("D27").Value = 105
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 419
$ws.Range("G27").Value = 423
$ws.Range("H27").Value = -0.945626477541
$ws.Range("I27").Value = 3363
$ws.Range("J27").Value = 3279
$ws.Range("K27").Value = 2.561756633119
$ws.Range("L27").Value = 12.174783188792

$ws.Range("C28").Value = 24
$ws.Range("D28").Value = 20
$ws.Range("E28").Value = 20
$ws.Range("F28").Value = 103
$ws.Range("G28").Value = 130
$ws.Range("H28").Value = -20.76923076923
$ws.Range("I28").Value = 774
$ws.Range("J28").Value = 1066
$ws.Range("K28").Value = -27.392120075046
$ws.Range("L28").Value = -35.5
$ws.Range("M28").Value = -32.75412684622
$ws.Range("N28").Value = -79.796397807361

$ws.Range("C29").Value = 22
$ws.Range("D29").Value = 20
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 91
$ws.Range("G29").Value = 111
$ws.Range("H29").Value = -18.018018018018
$ws.Range("I29").Value = 655
$ws.Range("J29").Value = 885
$ws.Range("K29").Value = -25.988700564971
$ws.Range("L29").Value = -34.760956175298
$ws.Range("M29").Value = -31.125131440588
$ws.Range("N29").Value = -80.992455020313

$ws.Range("C30").Value = 5
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = -50
$ws.Range("F30").Value = 24
$ws.Range("G30").Value = 56
$ws.Range("H30").Value = -57.142857142857
$ws.Range("I30").Value = 303
$ws.Range("J30").Value = 436
$ws.Range("K30").Value = -30.504587155963
$ws.Range("L30").Value = -14.164305949008

